# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.969.43"
$ws.Range("E2").Value = "  +2.56%  "

$ws.Range("D3").Value = "2.465.74"
$ws.Range("E3").Value = "  +4.46%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'566.17"
$ws.Range("E5").Value = "  +1.40%  "

$ws.Range("D6").Value = "'142.79"
$ws.Range("E6").Value = "  +6.95%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("D9").Value = "2.464.26"
$ws.Range("E9").Value = "  +4.55%  "

$ws.Range("E10").Value = "  +1.65%  "

$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("E14").Value = "  +8.72%  "

$ws.Range("D15").Value = "2.908.81"

$ws.Range("D16").Value = "62.839.79"
$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("E17").Value = "  +3.48%  "

$ws.Range("D18").Value = "2.468.01"
$ws.Range("E18").Value = "  +4.67%  "

$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  +3.57%  "

$ws.Range("D20").Value = "'340.30"
$ws.Range("E20").Value = "  +5.61%  "

$ws.Range("E21").Value = "  +2.35%  "

$ws.Range("D22").Value = "'6.81"
$ws.Range("E22").Value = "  +1.23%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "'65.47"
$ws.Range("E24").Value = "  +1.35%  "

$ws.Range("D25").Value = "'0.171"
$ws.Range("E25").Value = "  -1.10%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +3.24%  "

$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").Value = "'1.38"
$ws.Range("E29").Value = "  +6.69%  "

$ws.Range("D30").Value = "'6.81"
$ws.Range("E30").Value = "  +10.06%  "

$ws.Range("E31").Value = "  +4.66%  "

$ws.Range("D32").Value = "0.0₃0795"
$ws.Range("E32").Value = "  +6.10%  "

$ws.Range("D33").Value = "'175.08"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("E34").Value = "  +8.52%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("E36").Value = "  +3.08%  "

$ws.Range("D37").Value = "'373.95"
$ws.Range("E37").Value = "  +11.49%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").Value = "'4.37"
$ws.Range("E39").Value = "  +3.35%  "

$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +8.50%  "

$ws.Range("D42").Value = "'40.31"
$ws.Range("E42").Value = "  +5.16%  "

$ws.Range("D43").Value = "'150.00"
$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("E44").Value = "  +3.05%  "

$ws.Range("D45").Value = "'20.54"
$ws.Range("E45").Value = "  +4.51%  "

$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "'0.0517"
$ws.Range("E48").Value = "  +2.21%  "

$ws.Range("D49").Value = "'0.0227"
$ws.Range("E49").Value = "  +3.60%  "

$ws.Range("D50").Value = "0.0₆0233"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").Value = "'17.88"
$ws.Range("E51").Value = "  +2.51%  "
